$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 88.72291666666666
$ws.Cells.Item(2, 8).Value = 266.16875
$ws.Cells.Item(2, 9).Value = 0.7675060578750151
$ws.Cells.Item(2, 10).Value = 0.7675060578750152
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 63.76294533333333
$ws.Cells.Item(2, 14).Value = 191.288836
$ws.Cells.Item(2, 15).Value = 0.6446527016991613
$ws.Cells.Item(2, 16).Value = 0.6446527016991614
$ws.Cells.Item(2, 17).Value = 5657.234485230555
$ws.Cells.Item(2, 18).Value = 50915.110367075
$ws.Cells.Item(2, 19).Value = 0.4947748537796014
$ws.Cells.Item(2, 20).Value = 0.4947748537796016

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 88.72291666666666
$ws.Cells.Item(3, 8).Value = 266.16875
$ws.Cells.Item(3, 9).Value = 0.7675060578750151
$ws.Cells.Item(3, 10).Value = 0.7675060578750152
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 22.32219333333333
$ws.Cells.Item(3, 14).Value = 66.96658
$ws.Cells.Item(3, 15).Value = 0.2256806388876402
$ws.Cells.Item(3, 16).Value = 0.2256806388876402
$ws.Cells.Item(3, 17).Value = 1980.490098930555
$ws.Cells.Item(3, 18).Value = 17824.410890375
$ws.Cells.Item(3, 19).Value = 0.1732112574913675
$ws.Cells.Item(3, 20).Value = 0.1732112574913676

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 88.72291666666666
$ws.Cells.Item(4, 8).Value = 266.16875
$ws.Cells.Item(4, 9).Value = 0.7675060578750151
$ws.Cells.Item(4, 10).Value = 0.7675060578750152
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 12.825399
$ws.Cells.Item(4, 14).Value = 38.476197
$ws.Cells.Item(4, 15).Value = 0.1296666594131984
$ws.Cells.Item(4, 16).Value = 0.1296666594131984
$ws.Cells.Item(4, 17).Value = 1137.90680669375
$ws.Cells.Item(4, 18).Value = 10241.16126024375
$ws.Cells.Item(4, 19).Value = 0.09951994660404613
$ws.Cells.Item(4, 20).Value = 0.09951994660404614

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 17.91585
$ws.Cells.Item(5, 8).Value = 53.74755
$ws.Cells.Item(5, 9).Value = 0.1549827702197958
$ws.Cells.Item(5, 10).Value = 0.1549827702197958
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 63.76294533333333
$ws.Cells.Item(5, 14).Value = 191.288836
$ws.Cells.Item(5, 15).Value = 0.6446527016991613
$ws.Cells.Item(5, 16).Value = 0.6446527016991614
$ws.Cells.Item(5, 17).Value = 1142.3673641502
$ws.Cells.Item(5, 18).Value = 10281.3062773518
$ws.Cells.Item(5, 19).Value = 0.09991006153901169
$ws.Cells.Item(5, 20).Value = 0.09991006153901173

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 17.91585
$ws.Cells.Item(6, 8).Value = 53.74755
$ws.Cells.Item(6, 9).Value = 0.1549827702197958
$ws.Cells.Item(6, 10).Value = 0.1549827702197958
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 22.32219333333333
$ws.Cells.Item(6, 14).Value = 66.96658
$ws.Cells.Item(6, 15).Value = 0.2256806388876402
$ws.Cells.Item(6, 16).Value = 0.2256806388876402
$ws.Cells.Item(6, 17).Value = 399.921067431
$ws.Cells.Item(6, 18).Value = 3599.289606879
$ws.Cells.Item(6, 19).Value = 0.03497661059977985
$ws.Cells.Item(6, 20).Value = 0.03497661059977986

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 17.91585
$ws.Cells.Item(7, 8).Value = 53.74755
$ws.Cells.Item(7, 9).Value = 0.1549827702197958
$ws.Cells.Item(7, 10).Value = 0.1549827702197958
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 12.825399
$ws.Cells.Item(7, 14).Value = 38.476197
$ws.Cells.Item(7, 15).Value = 0.1296666594131984
$ws.Cells.Item(7, 16).Value = 0.1296666594131984
$ws.Cells.Item(7, 17).Value = 229.77792467415
$ws.Cells.Item(7, 18).Value = 2068.00132206735
$ws.Cells.Item(7, 19).Value = 0.02009609808100425
$ws.Cells.Item(7, 20).Value = 0.02009609808100426

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 8.960212333333333
$ws.Cells.Item(8, 8).Value = 26.880637
$ws.Cells.Item(8, 9).Value = 0.077511171905189
$ws.Cells.Item(8, 10).Value = 0.07751117190518901
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 63.76294533333333
$ws.Cells.Item(8, 14).Value = 191.288836
$ws.Cells.Item(8, 15).Value = 0.6446527016991613
$ws.Cells.Item(8, 16).Value = 0.6446527016991614
$ws.Cells.Item(8, 17).Value = 571.3295291853924
$ws.Cells.Item(8, 18).Value = 5141.965762668532
$ws.Cells.Item(8, 19).Value = 0.04996778638054822
$ws.Cells.Item(8, 20).Value = 0.04996778638054823

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 8.960212333333333
$ws.Cells.Item(9, 8).Value = 26.880637
$ws.Cells.Item(9, 9).Value = 0.077511171905189
$ws.Cells.Item(9, 10).Value = 0.07751117190518901
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 22.32219333333333
$ws.Cells.Item(9, 14).Value = 66.96658
$ws.Cells.Item(9, 15).Value = 0.2256806388876402
$ws.Cells.Item(9, 16).Value = 0.2256806388876402
$ws.Cells.Item(9, 17).Value = 200.0115920123844
$ws.Cells.Item(9, 18).Value = 1800.10432811146
$ws.Cells.Item(9, 19).Value = 0.01749277079649276
$ws.Cells.Item(9, 20).Value = 0.01749277079649276

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 8.960212333333333
$ws.Cells.Item(10, 8).Value = 26.880637
$ws.Cells.Item(10, 9).Value = 0.077511171905189
$ws.Cells.Item(10, 10).Value = 0.07751117190518901
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 12.825399
$ws.Cells.Item(10, 14).Value = 38.476197
$ws.Cells.Item(10, 15).Value = 0.1296666594131984
$ws.Cells.Item(10, 16).Value = 0.1296666594131984
$ws.Cells.Item(10, 17).Value = 114.918298299721
$ws.Cells.Item(10, 18).Value = 1034.264684697489
$ws.Cells.Item(10, 19).Value = 0.01005061472814802
$ws.Cells.Item(10, 20).Value = 0.01005061472814802
